$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New header column H: "Save" — copy the header formatting from the
# neighboring "sum" header (G1) so it picks up the same bold/bordered/
# centered style, then set its text.
$ws.Range("G1").Copy()
$ws.Range("H1").PasteSpecial(-4122)  # xlPasteFormats
$ws.Range("H1").Value = "Save"

# New data column H row 2: numeric value 1
$ws.Range("H2").Value = 1
